$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# praclen (column I) bumped from 4 to 5 for the existing training trials
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 5

# Append a new training trial row (row 6)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 61
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim1_2"

# Reset the cursor to A1 (was left on I5 in the prior save)
$ws.Range("A1").Select()
